$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap match data (columns F:V) between rows 79 and 80 ---
# Row 79 becomes the former row-80 match (Giresunspor vs Boluspor)
$ws.Range("F79").Value = "Giresunspor"
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = "Boluspor"
$ws.Range("I79").Value = 1
$ws.Range("J79").Value = 2.93
$ws.Range("K79").Value = "16/10/2023 20:42"
$ws.Range("L79").Value = 3.53
$ws.Range("M79").Value = "22/10/2023 14:59"
$ws.Range("N79").Value = 3.18
$ws.Range("O79").Value = "16/10/2023 20:42"
$ws.Range("P79").Value = 3.4
$ws.Range("Q79").Value = "22/10/2023 14:59"
$ws.Range("R79").Value = 2.52
$ws.Range("S79").Value = "16/10/2023 20:42"
$ws.Range("T79").Value = 2.14
$ws.Range("U79").Value = "22/10/2023 14:59"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/turkey/1-lig/giresunspor-boluspor/QPUJRWtI/"

# Row 80 becomes the former row-79 match (Umraniyespor vs Manisa FK)
$ws.Range("F80").Value = "Umraniyespor"
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = "Manisa FK"
$ws.Range("I80").Value = 2
$ws.Range("J80").Value = 2.76
$ws.Range("K80").Value = "16/10/2023 01:12"
$ws.Range("L80").Value = 3.09
$ws.Range("M80").Value = "22/10/2023 14:51"
$ws.Range("N80").Value = 3.29
$ws.Range("O80").Value = "16/10/2023 01:12"
$ws.Range("P80").Value = 3.23
$ws.Range("Q80").Value = "22/10/2023 14:51"
$ws.Range("R80").Value = 2.58
$ws.Range("S80").Value = "16/10/2023 01:12"
$ws.Range("T80").Value = 2.43
$ws.Range("U80").Value = "22/10/2023 14:51"
$ws.Range("V80").Value = "https://www.betexplorer.com/football/turkey/1-lig/umraniyespor-manisa-fk/lKZfXzDt/"

# --- Swap match data (columns F:V) between rows 106 and 107 ---
# Row 106 becomes the former row-107 match (Boluspor vs Sakaryaspor)
$ws.Range("F106").Value = "Boluspor"
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = "Sakaryaspor"
$ws.Range("I106").Value = 3
$ws.Range("J106").Value = 2.79
$ws.Range("K106").Value = "05/11/2023 14:12"
$ws.Range("L106").Value = 2.69
$ws.Range("M106").Value = "12/11/2023 11:01"
$ws.Range("N106").Value = 3.15
$ws.Range("O106").Value = "05/11/2023 14:12"
$ws.Range("P106").Value = 3.01
$ws.Range("Q106").Value = "12/11/2023 10:41"
$ws.Range("R106").Value = 2.65
$ws.Range("S106").Value = "05/11/2023 14:12"
$ws.Range("T106").Value = 2.95
$ws.Range("U106").Value = "12/11/2023 11:01"
$ws.Range("V106").Value = "https://www.betexplorer.com/football/turkey/1-lig/boluspor-sakaryaspor/IcM3n7Ip/"

# Row 107 becomes the former row-106 match (Manisa FK vs Bandirmaspor)
$ws.Range("F107").Value = "Manisa FK"
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = "Bandirmaspor"
$ws.Range("I107").Value = 2
$ws.Range("J107").Value = 2.18
$ws.Range("K107").Value = "05/11/2023 11:42"
$ws.Range("L107").Value = 2.54
$ws.Range("M107").Value = "12/11/2023 11:21"
$ws.Range("N107").Value = 3.41
$ws.Range("O107").Value = "05/11/2023 11:42"
$ws.Range("P107").Value = 3.35
$ws.Range("Q107").Value = "12/11/2023 11:25"
$ws.Range("R107").Value = 3.31
$ws.Range("S107").Value = "05/11/2023 11:42"
$ws.Range("T107").Value = 2.84
$ws.Range("U107").Value = "12/11/2023 11:25"
$ws.Range("V107").Value = "https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-bandirmaspor/Eq666S9G/"

# --- Append 3 new match rows (121, 122, 123) ---
# Seed the new rows from row 120 so they inherit the same cell styles
# (bold+bordered index column, date-formatted match-date column), then
# overwrite every cell with the real values.
$ws.Range("A120:V120").Copy($ws.Range("A121"))
$ws.Range("A120:V120").Copy($ws.Range("A122"))
$ws.Range("A120:V120").Copy($ws.Range("A123"))

# Row 121: Keciorengucu 1 - 3 Sakaryaspor
$ws.Range("A121").Value = 120
$ws.Range("B121").Value = "turkey"
$ws.Range("C121").Value = "1-lig"
$ws.Range("D121").Value = "2023-2024"
$ws.Range("E121").Value = 45262.47916666666
$ws.Range("F121").Value = "Keciorengucu"
$ws.Range("G121").Value = 1
$ws.Range("H121").Value = "Sakaryaspor"
$ws.Range("I121").Value = 3
$ws.Range("J121").Value = 3.77
$ws.Range("K121").Value = "27/11/2023 18:12"
$ws.Range("L121").Value = 2.74
$ws.Range("M121").Value = "02/12/2023 11:21"
$ws.Range("N121").Value = 3.54
$ws.Range("O121").Value = "27/11/2023 18:12"
$ws.Range("P121").Value = 3.31
$ws.Range("Q121").Value = "02/12/2023 11:21"
$ws.Range("R121").Value = 1.97
$ws.Range("S121").Value = "27/11/2023 18:12"
$ws.Range("T121").Value = 2.65
$ws.Range("U121").Value = "02/12/2023 11:21"
$ws.Range("V121").Value = "https://www.betexplorer.com/football/turkey/1-lig/keciorengucu-sakaryaspor/Y7aGqQJ0/"

# Row 122: Tuzlaspor 1 - 1 Bandirmaspor
$ws.Range("A122").Value = 121
$ws.Range("B122").Value = "turkey"
$ws.Range("C122").Value = "1-lig"
$ws.Range("D122").Value = "2023-2024"
$ws.Range("E122").Value = 45262.47916666666
$ws.Range("F122").Value = "Tuzlaspor"
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = "Bandirmaspor"
$ws.Range("I122").Value = 1
$ws.Range("J122").Value = 3.95
$ws.Range("K122").Value = "26/11/2023 17:12"
$ws.Range("L122").Value = 4.81
$ws.Range("M122").Value = "02/12/2023 11:25"
$ws.Range("N122").Value = 3.62
$ws.Range("O122").Value = "26/11/2023 17:12"
$ws.Range("P122").Value = 3.75
$ws.Range("Q122").Value = "02/12/2023 11:25"
$ws.Range("R122").Value = 1.9
$ws.Range("S122").Value = "26/11/2023 17:12"
$ws.Range("T122").Value = 1.75
$ws.Range("U122").Value = "02/12/2023 11:25"
$ws.Range("V122").Value = "https://www.betexplorer.com/football/turkey/1-lig/tuzlaspor-bandirmaspor/fNmPs4lD/"

# Row 123: Boluspor 0 - 3 Eyupspor
$ws.Range("A123").Value = 122
$ws.Range("B123").Value = "turkey"
$ws.Range("C123").Value = "1-lig"
$ws.Range("D123").Value = "2023-2024"
$ws.Range("E123").Value = 45262.58333333334
$ws.Range("F123").Value = "Boluspor"
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = "Eyupspor"
$ws.Range("I123").Value = 3
$ws.Range("J123").Value = 5.18
$ws.Range("K123").Value = "26/11/2023 14:13"
$ws.Range("L123").Value = 4.71
$ws.Range("M123").Value = "02/12/2023 13:55"
$ws.Range("N123").Value = 3.68
$ws.Range("O123").Value = "26/11/2023 14:13"
$ws.Range("P123").Value = 3.63
$ws.Range("Q123").Value = "02/12/2023 13:55"
$ws.Range("R123").Value = 1.68
$ws.Range("S123").Value = "26/11/2023 14:13"
$ws.Range("T123").Value = 1.79
$ws.Range("U123").Value = "02/12/2023 13:55"
$ws.Range("V123").Value = "https://www.betexplorer.com/football/turkey/1-lig/boluspor-eyupspor/K4eCp64f/"
